$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)
$p1.Style = "Compact"

for ($i = 0; $i -lt 4; $i++) {
    $p1.Range.InsertParagraphAfter()
}
for ($i = 2; $i -le 5; $i++) {
    $d.Paragraphs.Item($i).Style = "Normal"
}

$p5 = $d.Paragraphs.Item(5)
$r = $p5.Range.Duplicate()
$r.Collapse(1)
Write-Host ("r start/end right before Add: " + $r.Start + " " + $r.End)
$d.Bookmarks.Add("_GoBack", $r)
Write-Host ("r start/end right after Add: " + $r.Start + " " + $r.End)

# check bookmark's own start/end via iterating the (hidden) bookmark
$d.Bookmarks.ShowHidden = $true
Write-Host ("bookmarks count: " + $d.Bookmarks.Count)
